# Apply the "update samples and readme.md" commit:
# adds two new sheets ("Dungeons" and "Npcs") with sample BakingSheet data,
# and updates the selection state on the existing sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1) Add "Dungeons" sheet after "Monsters" (becomes the 4th sheet)
# ---------------------------------------------------------------
$afterSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$dungeons = $wb.Worksheets.Add($null, $afterSheet)
$dungeons.Name = "Dungeons"

$dungeons.Cells.Item(1,1).Value = "Id"
$dungeons.Cells.Item(1,2).Value = "Name"
$dungeons.Cells.Item(1,3).Value = "Monsters:1"
$dungeons.Cells.Item(1,4).Value = "Monsters:2"
$dungeons.Cells.Item(1,5).Value = "Monsters:3"
$dungeons.Cells.Item(1,6).Value = "Items:1"
$dungeons.Cells.Item(1,7).Value = "Items:2"

$dungeons.Cells.Item(2,1).Value = "DUNGEON001"
$dungeons.Cells.Item(2,2).Value = "Easy Field"
$dungeons.Cells.Item(2,3).Value = "MONSTER001"
$dungeons.Cells.Item(2,6).Value = "ITEM_POTION001"
$dungeons.Cells.Item(2,7).Value = "ITEM_LVUP001"

$dungeons.Cells.Item(3,1).Value = "DUNGEON002"
$dungeons.Cells.Item(3,2).Value = "Expert Zone"
$dungeons.Cells.Item(3,3).Value = "MONSTER001"
$dungeons.Cells.Item(3,4).Value = "MONSTER002"
$dungeons.Cells.Item(3,6).Value = "ITEM_POTION002"
$dungeons.Cells.Item(3,7).Value = "ITEM_LVUP002"

$dungeons.Cells.Item(4,1).Value = "DUNGEON003"
$dungeons.Cells.Item(4,2).Value = "Dragon" + [char]0x2019 + "s Nest"
$dungeons.Cells.Item(4,3).Value = "MONSTER003"
$dungeons.Cells.Item(4,4).Value = "MONSTER004"
$dungeons.Cells.Item(4,5).Value = "MONSTER005"
$dungeons.Cells.Item(4,6).Value = "ITEM_LVUP003"

# ---------------------------------------------------------------
# 2) Add "Npcs" sheet after "Dungeons" (becomes the 5th sheet)
# ---------------------------------------------------------------
$npcs = $wb.Worksheets.Add($null, $dungeons)
$npcs.Name = "Npcs"

$npcs.Cells.Item(1,1).Value = "Id"
$npcs.Cells.Item(1,2).Value = "Name"
$npcs.Cells.Item(1,3).Value = "Texts:Greeting"
$npcs.Cells.Item(1,4).Value = "Texts:Purchasing"
$npcs.Cells.Item(1,5).Value = "Texts:Leaving"

$npcs.Cells.Item(2,1).Value = "NPC001"
$npcs.Cells.Item(2,2).Value = "Fat Baker"
$npcs.Cells.Item(2,3).Value = "Morning traveler!"
$npcs.Cells.Item(2,4).Value = "Thank you!"
$npcs.Cells.Item(2,5).Value = "Come again!"

$npcs.Cells.Item(3,1).Value = "NPC002"
$npcs.Cells.Item(3,2).Value = "Blacksmith"
$npcs.Cells.Item(3,3).Value = "G" + [char]0x2019 + "day!"
$npcs.Cells.Item(3,4).Value = "Good choice."
$npcs.Cells.Item(3,5).Value = "Take care."

$npcs.Cells.Item(4,1).Value = "NPC003"
$npcs.Cells.Item(4,2).Value = "Potion Master"
$npcs.Cells.Item(4,3).Value = "What do you want?"
$npcs.Cells.Item(4,4).Value = "Take it already."
$npcs.Cells.Item(4,5).Value = "Don" + [char]0x2019 + "t come again."

# ---------------------------------------------------------------
# 3) Update selection state on each sheet (also drives activeTab /
#    tabSelected bookkeeping since the last Select() wins).
# ---------------------------------------------------------------
$items = $wb.Worksheets.Item("Items")
[void]$items.Range("A2").Select()

$monsters = $wb.Worksheets.Item("Monsters")
[void]$monsters.Range("A2").Select()

[void]$dungeons.Range("E29").Select()

[void]$npcs.Range("E2").Select()
